$p = $ppt.ActivePresentation

# --- 1. Update the hard-coded "datetimeFigureOut" field text on the Slide
#        Master's Date placeholder (format: 2010/5/19 -> 2010/8/9) ---
$masterDateShape = $p.SlideMaster.Shapes.Item(3)
$masterDateShape.TextFrame.TextRange.Text = "2010/8/9"

# --- 2. Update the same field on all 11 Slide Layouts ---
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "2010/5/19") {
                $shp.TextFrame.TextRange.Text = "2010/8/9"
            }
        }
    }
}

# --- 3. Update the Notes Master's date field (format: 5/19/2010 -> 8/9/2010) ---
$notesMaster = $p.NotesMaster
for ($j = 1; $j -le $notesMaster.Shapes.Count; $j++) {
    $shp = $notesMaster.Shapes.Item($j)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "5/19/2010") {
            $shp.TextFrame.TextRange.Text = "8/9/2010"
        }
    }
}

# --- 4. Slide 3: split the "Duration [start_time, end_time[ (secs)" table
#        header cell into two paragraphs and fix up the wording:
#          "Duration [start_time, end_time[ (secs)"
#        becomes
#          "Duration "                     (paragraph 1)
#          "[start-time, end-time[ (secs)" (paragraph 2)
$slide3 = $p.Slides.Item(3)
$tableShape = $slide3.Shapes.Item(4)
$cell = $tableShape.Table.Cell(1, 1)
$cell.Shape.TextFrame.TextRange.Text = "Duration `r[start-time, end-time[ (secs)"
